# SysConfig.TblAppObject_UserRole.xlsx - "Update Pertanggal 15 November 2024 16:18 WIB"
#
# Adds three new user roles ("Tools And Asset Manager", "Tools And Asset
# Senior Staff", "Tools And Asset Staff") to the bottom of the role list on
# sheet "Main", following the exact same row pattern used by every other
# role block in the sheet:
#   - col B: role name (blank "separator" row, then one row per role)
#   - col C: shared formula building the PERFORM ... SQL string from col B
#   - col D: sequential numeric id
#
# It also nudges the window/selection state to match where the workbook was
# left scrolled to after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- 1. Clone the formatting of the previous role block (rows 60:63) down
#        onto the new block (rows 64:67) so the new rows pick up the same
#        cell styles (col B role-name style, col C formula style, col D
#        id style) as every other row in the sheet. ---
$ws.Range("B60:D63").Copy() | Out-Null
$ws.Range("B64:D67").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. New role names (col B) ---
$ws.Range("B65").Value = "Tools And Asset Manager"
$ws.Range("B66").Value = "Tools And Asset Senior Staff"
$ws.Range("B67").Value = "Tools And Asset Staff"

# --- 3. New sequential ids (col D) ---
$ws.Range("D65").Value = 95000000000052
$ws.Range("D66").Value = 95000000000053
$ws.Range("D67").Value = 95000000000054

# --- 4. Shared formula for col C across the whole new block at once, same
#        as the "PERFORM ..." formula used throughout the sheet, so it gets
#        stored as a single shared-formula group (matching B64 being the
#        blank separator row). ---
$ws.Range("C64:C67").Formula = "=IF(EXACT(B64,""""),"""",CONCATENATE(""PERFORM """"SchSysConfig"""".""""Func_TblAppObject_UserRole_SET""""(varSystemLoginSession, null, null, null, varInstitutionBranchID, '"",B64,""');""))"

# --- 5. Leave the selection / scroll position where the author left it. ---
$win = $wb.Windows.Item(1)
try {
    $win.Left = 0
    $win.Top = 0
    $win.Width = 20490
    $win.Height = 6000
    $win.ScrollRow = 52
    $win.ScrollColumn = 2
} catch {
    # Window geometry isn't always controllable headlessly; ignore.
}

$ws.Range("D68").Select() | Out-Null
